# Atualização de bases das ligas, do dia: 12-04-2024 às 20:28
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The match row (id=144, matchId 6769308, NK Rudes v Slaven Belupo) was
# removed from the source feed. Delete its entire worksheet row; this
# shifts every row below it up by one (old rows 147-150 become 146-149).
$ws.Rows(146).EntireRow.Delete()

# Column A is a sequential row id (row number - 2). Restore it for every
# row that shifted up.
$ws.Range("A146").Value = 144
$ws.Range("A147").Value = 145
$ws.Range("A148").Value = 146
$ws.Range("A149").Value = 147

# Refreshed odds for the (now) row 146 - NK Varazdin v NK Lokomotiva Zagreb
$ws.Range("N146").Value = 2.75
$ws.Range("O146").Value = 3.3
$ws.Range("R146").Value = 1.975
$ws.Range("S146").Value = 2.5
$ws.Range("U146").Value = 1.875
$ws.Range("V146").Value = 1.975

# Refreshed odds for the (now) row 147 - HNK Gorica v Dinamo Zagreb
$ws.Range("N147").Value = 10
$ws.Range("O147").Value = 5
$ws.Range("P147").Value = 1.285
$ws.Range("R147").Value = 1.875
$ws.Range("S147").Value = 1.975

# Row 148 (Hajduk Split v NK Osijek) only needed the id fix above.

# Refreshed odds for the (now) row 149 - Istra 1961 v HNK Rijeka
$ws.Range("N149").Value = 6
$ws.Range("O149").Value = 3.75
$ws.Range("P149").Value = 1.55
$ws.Range("Q149").Value = 1
$ws.Range("R149").Value = 1.825
$ws.Range("S149").Value = 2.025
